$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the value for April 2025 (row 125, column B)
$ws.Range("B125").NumberFormat = "@"
$ws.Range("B125").Value = "65,227"
$ws.Range("B125").Style = "Normal"

# Append the new row for May 2025 (row 126)
$ws.Range("A126").NumberFormat = "@"
$ws.Range("A126").Value = "2025 MAY"
$ws.Range("A126").Style = "Normal"

$ws.Range("B126").NumberFormat = "@"
$ws.Range("B126").Value = "65,065"
$ws.Range("B126").Style = "Normal"
